# FormatoExcel2007.xlsx edits:
#   - "Se minimizan las opciones a diligenciar en la plantilla de
#      configuracion": the "Tipo de servicio" / "Se guardara en" header +
#      value cells are no longer needed on the "Paises" and "Clientes"
#      sheets, so their B/C cells in rows 1-2 are cleared out entirely.
#   - "Cuando el usuario tiene una sesion abierta, se habilita la pagina
#      de inicio": "Paises" becomes the active sheet/tab (instead of
#      "Clientes"), with updated cursor positions on both sheets.

$wb = $excel.ActiveWorkbook

$wsPaises = $wb.Worksheets.Item("Paises")
$wsClientes = $wb.Worksheets.Item("Clientes")

# --- Minimize the configuration template on both sheets ---
$wsPaises.Range("B1").Clear() | Out-Null
$wsPaises.Range("C1").Clear() | Out-Null
$wsPaises.Range("B2").Clear() | Out-Null
$wsPaises.Range("C2").Clear() | Out-Null

$wsClientes.Range("B1").Clear() | Out-Null
$wsClientes.Range("C1").Clear() | Out-Null
$wsClientes.Range("B2").Clear() | Out-Null
$wsClientes.Range("C2").Clear() | Out-Null

# --- Update cursor position on "Clientes" (no longer the active tab) ---
$wsClientes.Range("A14").Select() | Out-Null

# --- Make "Paises" the active/home sheet with its own cursor position ---
$wsPaises.Activate() | Out-Null
$wsPaises.Range("C13").Select() | Out-Null
